$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "서연"
$ws.Range("E1").Value = "현빈"

$ws.Range("D2").Value = "현빈 병국"
$ws.Range("E2").Value = ""

$ws.Range("C3").Value = "유진 재현"
$ws.Range("D3").Value = "현빈"
$ws.Range("E3").Value = ""

$ws.Range("D4").Value = "준범"
$ws.Range("E4").Value = ""

$ws.Range("C5").Value = "서연 태훈"
$ws.Range("D5").Value = "준범"
$ws.Range("E5").Value = ""

$ws.Range("D6").Value = "준범 현빈"
$ws.Range("E6").Value = ""

$ws.Range("A7").Value = "한솔 희지"
$ws.Range("B7").Value = "한솔 희지"
$ws.Range("C7").Value = "준범"
$ws.Range("D7").Value = ""
$ws.Range("E7").Value = ""
